# This script reproduces the reordering of the "Boiler turbogenerator" /
# "Chilled water package" column groups and the "Stream-*" price columns,
# plus the updated Monte-Carlo result values, in the Monte Carlo sugarcane
# results workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 1 (section headers)
# ---------------------------------------------------------------------
# "Boiler turbogenerator" / "Chilled water package" groups swap places
# and change extent (AJ1:AL1 + AM1:AY1  ->  AJ1:AV1 + AW1:AY1).
$ws.Range("AJ1:AL1").UnMerge()
$ws.Range("AM1:AY1").UnMerge()

$ws.Range('AJ1').Value = 'Boiler turbogenerator'
$ws.Range("AM1").Value = ""
$ws.Range('AW1').Value = 'Chilled water package'

$ws.Range("AJ1:AV1").Merge()
$ws.Range("AW1:AY1").Merge()

# Stream price columns get reordered: three new streams (Stream-s160,
# Stream-s161, Stream-cooling tower chemicals) are inserted before
# Stream-makeup water, pushing the rest of the columns over (Stream-sugarcane
# moves later in the run, next to Stream-denaturant).
$ws.Range('BN1').Value = 'Stream-s160'
$ws.Range('BO1').Value = 'Stream-s161'
$ws.Range('BP1').Value = 'Stream-cooling tower chemicals'
$ws.Range('BQ1').Value = 'Stream-makeup water'
$ws.Range('BR1').Value = 'Stream-denaturant'
$ws.Range('BS1').Value = 'Stream-sugarcane'
$ws.Range('BT1').Value = 'Stream-enzyme'
$ws.Range('BU1').Value = 'Stream-lime'

# ---------------------------------------------------------------------
# Row 2 (variable headers) - follow the same reordering of the
# Boiler turbogenerator / Chilled water package fields.
# ---------------------------------------------------------------------
$ws.Range('AJ2').Value = 'Baghouse bags base cost [USD]'
$ws.Range('AK2').Value = 'Baghouse bags exponent'
$ws.Range('AL2').Value = 'Boiler base cost [USD]'
$ws.Range('AM2').Value = 'Boiler exponent'
$ws.Range('AN2').Value = 'Boiler electricity rate [kW / kg/hr]'
$ws.Range('AO2').Value = 'Deaerator base cost [USD]'
$ws.Range('AP2').Value = 'Deaerator exponent'
$ws.Range('AQ2').Value = 'Amine addition pkg base cost [USD]'
$ws.Range('AR2').Value = 'Amine addition pkg exponent'
$ws.Range('AS2').Value = 'Hot process water softener system base cost [USD]'
$ws.Range('AT2').Value = 'Hot process water softener system exponent'
$ws.Range('AU2').Value = 'Turbogenerator base cost [USD]'
$ws.Range('AV2').Value = 'Turbogenerator exponent'
$ws.Range('AW2').Value = 'Base cost [USD]'
$ws.Range('AX2').Value = 'Exponent'
$ws.Range('AY2').Value = 'Electricity rate [kW / kJ/hr]'

# ---------------------------------------------------------------------
# Row 4 (data values) - values travel together with their (now moved)
# column headers, and a handful of downstream results are recalculated.
# ---------------------------------------------------------------------
$ws.Range("AJ4").Value = 106.8491863396745
$ws.Range("AK4").Value = 1
$ws.Range("AL4").Value = 28550000
$ws.Range("AM4").Value = 0.6
$ws.Range("AN4").Value = 0.005743948115934742
$ws.Range("AO4").Value = 305000
$ws.Range("AQ4").Value = 40000
$ws.Range("AR4").Value = 0.6
$ws.Range("AS4").Value = 78000
$ws.Range("AT4").Value = 0.6
$ws.Range("AU4").Value = 9500000
$ws.Range("AV4").Value = 0.6
$ws.Range("AW4").Value = 1375000
$ws.Range("AX4").Value = 0.7
$ws.Range("AY4").Value = [double]"-4.328359737776564e-05"

$ws.Range("BN4").Value = 0.199375046806894
$ws.Range("BO4").Value = 4.995862254032183
$ws.Range("BP4").Value = 3
$ws.Range("BQ4").Value = 0.000254
$ws.Range("BR4").Value = 0.756
$ws.Range("BS4").Value = 0.03455
$ws.Range("BT4").Value = 0.5
$ws.Range("BU4").Value = 0.077

$ws.Range("CF4").Value = 0.1177528682799202
$ws.Range("CG4").Value = 65657604.31377867
$ws.Range("CH4").Value = 174443300.1044682
$ws.Range("CJ4").Value = 720922.9237370482
$ws.Range("CK4").Value = 48512.43283670145
$ws.Range("CL4").Value = 183439.3057284494
